# ===========================================================================
# Generate Report for Handback
#
# The localized files (zh-cn + de-de) have now been handed back in sync
# with en-US: every "Ready for handoff" status becomes "Handed back: in
# sync with en-US", and each language sheet records the localized target
# file, the handback xlf file name, and the handback timestamp.
# ===========================================================================

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status update: handback completed, all rows now in sync with en-US.
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# ---------------------------------------------------------------------------
# 2. Widen columns that now hold the longer status text / target-file links.
# ---------------------------------------------------------------------------
$wsOverview.Columns("E").ColumnWidth = 29.166666666666668
$wsOverview.Columns("F").ColumnWidth = 29.166666666666668

$wsZhCn.Columns("C").ColumnWidth = 29.166666666666668
$wsZhCn.Columns("I").ColumnWidth = 39.166666666666664
$wsZhCn.Columns("J").ColumnWidth = 39.166666666666664

$wsDeDe.Columns("C").ColumnWidth = 29.166666666666668
$wsDeDe.Columns("I").ColumnWidth = 39.166666666666664
$wsDeDe.Columns("J").ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# 3. zh-cn handback: Latest Target File / Latest Handback File / DateTime.
# ---------------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cacdb166865b23e85868a85cb7be2548a36fb23/e2e/1aba61aa-d259-415f-9447-4f8db1ef9e15.md", `
    "", "", "1aba61aa-d259-415f-9447-4f8db1ef9e15.md")
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("J2").Value = "1aba61aa-d259-415f-9447-4f8db1ef9e15.548d2abe1ab53c22add390a263f740f8913a8fea.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-29 03:01:32"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cacdb166865b23e85868a85cb7be2548a36fb23/e2e/1fe2c9e0-7f52-4499-a4c6-0e4c1466dce6.md", `
    "", "", "1fe2c9e0-7f52-4499-a4c6-0e4c1466dce6.md")
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Color = 15570276
$wsZhCn.Range("J3").Value = "1fe2c9e0-7f52-4499-a4c6-0e4c1466dce6.0537f948374ccd930d7dfab2b0d917ab8642c0b5.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-29 03:01:32"

# ---------------------------------------------------------------------------
# 4. de-de handback: Latest Target File / Latest Handback File / DateTime.
# ---------------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cacdb166865b23e85868a85cb7be2548a36fb23/e2e/1aba61aa-d259-415f-9447-4f8db1ef9e15.md", `
    "", "", "1aba61aa-d259-415f-9447-4f8db1ef9e15.md")
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("J2").Value = "1aba61aa-d259-415f-9447-4f8db1ef9e15.548d2abe1ab53c22add390a263f740f8913a8fea.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-29 03:01:39"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cacdb166865b23e85868a85cb7be2548a36fb23/e2e/1fe2c9e0-7f52-4499-a4c6-0e4c1466dce6.md", `
    "", "", "1fe2c9e0-7f52-4499-a4c6-0e4c1466dce6.md")
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Color = 15570276
$wsDeDe.Range("J3").Value = "1fe2c9e0-7f52-4499-a4c6-0e4c1466dce6.0537f948374ccd930d7dfab2b0d917ab8642c0b5.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-29 03:01:39"
